## Add Rothman state law link
##
## Inserts a new "Compact"-styled, numbered (numId=1003) list paragraph
## containing a hyperlink to "Rothman's Roadmap to the Right of Publicity",
## placed as the last bullet in the list of Matthews v. Wozencraft
## reference links (right after the "Texas Right of Publicity Law, citing
## Matthews." bullet, and right before the "JERRY E. SMITH, Circuit Judge:"
## paragraph that starts the opinion body).

$d = $word.ActiveDocument

# Locate the paragraph that should immediately follow the new bullet.
$anchor = $d.Content
$found = $anchor.Find.Execute("JERRY E. SMITH, Circuit Judge:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph 'JERRY E. SMITH, Circuit Judge:'"
}

# Insert the hyperlink text right before that paragraph (still inside the
# previous -- "Texas Right of Publicity Law, citing Matthews." -- paragraph
# for now; we'll split it into its own paragraph next).
$insertionPoint = $d.Range($anchor.Start, $anchor.Start)
$d.Hyperlinks.Add($insertionPoint, "http://www.rightofpublicityroadmap.com/", "", "", "Rothman’s Roadmap to the Right of Publicity") | Out-Null

# Re-find the "JERRY E. SMITH" paragraph (its start moved forward because
# text was just inserted before it) and split it into its own paragraph by
# inserting a paragraph break right before it.
$jerry = $d.Content
$jerry.Find.Execute("JERRY E. SMITH, Circuit Judge:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$breakPoint = $d.Range($jerry.Start, $jerry.Start)
$breakPoint.InsertParagraphBefore()

# Grab the paragraph that now holds just the new hyperlink (the paragraph
# right before "JERRY E. SMITH ...").
$jerry2 = $d.Content
$jerry2.Find.Execute("JERRY E. SMITH, Circuit Judge:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$jerryPara = $jerry2.Paragraphs(1)
$newPara = $jerryPara.Previous()

# Style it like the other reference bullets above it.
$newPara.Range.Style = "Compact"

# Reuse the same list (numId=1003) that the sibling bullets above use by
# copying the list template from one of them and continuing that same list.
$siblings = $d.Content
$siblings.Find.Execute("case on Westlaw", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$siblingPara = $siblings.Paragraphs(1)
$listTemplate = $siblingPara.Range.ListFormat.ListTemplate
$newPara.Range.ListFormat.ApplyListTemplate($listTemplate, $true)

Write-Output "New paragraph text: $($newPara.Range.Text)"
